# Updates the "Estado de Cuenta" (account statement) worksheet: the previous
# overdue periods (EC) are removed and replaced with the new set of periods,
# i.e. the "Periodo Mora" (column E) / "Valor Mora" (column F) pairs for the
# detail rows (16-30) are refreshed with the current database extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Periodo Mora, Valor Mora)
$ws.Range("E16").Value = "1804"
$ws.Range("F16").Value = 31249

$ws.Range("E17").Value = "1803"
$ws.Range("F17").Value = 31249

$ws.Range("E18").Value = "1907"
$ws.Range("F18").Value = 20979

$ws.Range("E19").Value = "1907"
$ws.Range("F19").Value = 1104

$ws.Range("E20").Value = "1906"
$ws.Range("F20").Value = 33125

$ws.Range("E21").Value = "1905"
$ws.Range("F21").Value = 31249

$ws.Range("E22").Value = "1904"
$ws.Range("F22").Value = 31249

$ws.Range("E23").Value = "1903"
$ws.Range("F23").Value = 31249

$ws.Range("E24").Value = "1902"
$ws.Range("F24").Value = 31249

$ws.Range("E25").Value = "1901"
$ws.Range("F25").Value = 31249

$ws.Range("E26").Value = "1812"
$ws.Range("F26").Value = 31249

$ws.Range("E27").Value = "1811"
$ws.Range("F27").Value = 31249

$ws.Range("E28").Value = "1810"
$ws.Range("F28").Value = 31249

$ws.Range("E29").Value = "1809"
$ws.Range("F29").Value = 31249

$ws.Range("E30").Value = "1808"
$ws.Range("F30").Value = 31249
